$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 data (20 minute trade that didn't need a close price / yahoo lookup)
$ws.Range("A4").Value = 10017
$ws.Range("B4").Value = 10002
$ws.Range("C4").Value = 80.45
$ws.Range("D4").Value = 80.569999999999993
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.15

$ws.Range("G4").Value = 42608.624062499999
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"

$ws.Range("H4").Value = $true
